# Add data for 2021-10-31
# Updates the "Through 2021-10-22" workbook to "Through 2021-10-23",
# incrementing counts for the new day's carjacking data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet and update the header label / shared string text
$ws.Name = "Through 2021-10-23"
$ws.Range("B1").Value = "October 2021 (through October 23)"

# Apply the updated cell values (existing and newly populated cells)
$ws.Range("L2").Value = 16
$ws.Range("V2").Value = 5
$ws.Range("B3").Value = 12
$ws.Range("B4").Value = 13
$ws.Range("AP6").Value = 2
$ws.Range("B9").Value = 7
$ws.Range("L9").Value = 2
$ws.Range("V9").Value = 4
$ws.Range("BJ13").Value = 3
$ws.Range("B18").Value = 2
$ws.Range("L22").Value = 2
$ws.Range("AP24").Value = 1
$ws.Range("B25").Value = 1
$ws.Range("B26").Value = 1
$ws.Range("B27").Value = 4
$ws.Range("B33").Value = 2
$ws.Range("B41").Value = 4
$ws.Range("V41").Value = 3
$ws.Range("L42").Value = 1
$ws.Range("L47").Value = 1
$ws.Range("L66").Value = 4
$ws.Range("L67").Value = 1
$ws.Range("AP67").Value = 2
$ws.Range("B71").Value = 1
$ws.Range("AF71").Value = 1
$ws.Range("L87").Value = 1
$ws.Range("AZ92").Value = 1
$ws.Range("AZ94").Value = 1
